$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 20.66320133773953
$ws.Range("C2").Value = 13.14330772211953
$ws.Range("D2").Value = 6.290061924380096
$ws.Range("F2").Value = 55.68467830992864
$ws.Range("G2").Value = 3.718094583745748
$ws.Range("J2").Value = 11.53462253742862
$ws.Range("B3").Value = 20.30961541451096
$ws.Range("C3").Value = 12.82187629189725
$ws.Range("D3").Value = 6.344075758192185
$ws.Range("F3").Value = 54.74634869498082
$ws.Range("G3").Value = 3.723363054668166
$ws.Range("J3").Value = 11.5279750397931
$ws.Range("B4").Value = 20.09857968690556
$ws.Range("C4").Value = 12.62734140610034
$ws.Range("D4").Value = 6.378993313946221
$ws.Range("F4").Value = 54.1736638812817
$ws.Range("G4").Value = 3.726757141330943
$ws.Range("J4").Value = 11.52653860730808
$ws.Range("B5").Value = 20.01423958151413
$ws.Range("C5").Value = 12.54892781398593
$ws.Range("D5").Value = 6.393664511047524
$ws.Range("F5").Value = 53.94139373803116
$ws.Range("G5").Value = 3.728180492328934
$ws.Range("J5").Value = 11.52661723180128
$ws.Range("B6").Value = 20.00033896975118
$ws.Range("C6").Value = 12.53596372728102
$ws.Range("D6").Value = 6.396127385646406
$ws.Range("F6").Value = 53.90289895017241
$ws.Range("G6").Value = 3.728419274361172
$ws.Range("J6").Value = 11.52667034860521
$ws.Range("B7").Value = 20.09743535758361
$ws.Range("C7").Value = 12.62628020679212
$ws.Range("D7").Value = 6.379189383464745
$ws.Range("F7").Value = 54.17052662701489
$ws.Range("G7").Value = 3.726776173978275
$ws.Range("J7").Value = 11.52653698103186
$ws.Range("B8").Value = 20.54010628634994
$ws.Range("C8").Value = 13.03198137613031
$ws.Range("D8").Value = 6.308322829995163
$ws.Range("F8").Value = 55.36056691581162
$ws.Range("G8").Value = 3.719878223307663
$ws.Range("J8").Value = 11.53178088827188
$ws.Range("B9").Value = 21.45036472302838
$ws.Range("C9").Value = 13.84351431699975
$ws.Range("D9").Value = 6.183206894766672
$ws.Range("F9").Value = 57.71130300559327
$ws.Range("G9").Value = 3.707605819909949
$ws.Range("J9").Value = 11.56308947925237
$ws.Range("B10").Value = 22.13682572235752
$ws.Range("C10").Value = 14.44116687294514
$ws.Range("D10").Value = 6.099655659630089
$ws.Range("F10").Value = 59.4351234678487
$ws.Range("G10").Value = 3.699341515623832
$ws.Range("J10").Value = 11.59893887745762
$ws.Range("B11").Value = 22.45128466010679
$ws.Range("C11").Value = 14.71173917510376
$ws.Range("D11").Value = 6.063449206169252
$ws.Range("F11").Value = 60.21581447316764
$ws.Range("G11").Value = 3.695742508613618
$ws.Range("J11").Value = 11.61803533000121
$ws.Range("B12").Value = 22.57054149884845
$ws.Range("C12").Value = 14.81388654878078
$ws.Range("D12").Value = 6.049996809653449
$ws.Range("F12").Value = 60.51071480823885
$ws.Range("G12").Value = 3.694402521627345
$ws.Range("J12").Value = 11.62566668110022
$ws.Range("B13").Value = 22.54485159835963
$ws.Range("C13").Value = 14.79190307843675
$ws.Range("D13").Value = 6.05288255292436
$ws.Range("F13").Value = 60.44723878938277
$ws.Range("G13").Value = 3.694690097383893
$ws.Range("J13").Value = 11.62400536645519
$ws.Range("B14").Value = 22.46109318494256
$ws.Range("C14").Value = 14.7201498483498
$ws.Range("D14").Value = 6.062337299226948
$ws.Range("F14").Value = 60.2400920227613
$ws.Range("G14").Value = 3.695631809676757
$ws.Range("J14").Value = 11.61865515782564
$ws.Range("B15").Value = 22.40980797387185
$ws.Range("C15").Value = 14.67615466828879
$ws.Range("D15").Value = 6.068162206095689
$ws.Range("F15").Value = 60.11310686547871
$ws.Range("G15").Value = 3.696211609320422
$ws.Range("J15").Value = 11.61543004688422
$ws.Range("B16").Value = 22.11630730712486
$ws.Range("C16").Value = 14.42344721842276
$ws.Range("D16").Value = 6.102058012204572
$ws.Range("F16").Value = 59.38401356436166
$ws.Range("G16").Value = 3.699579931450255
$ws.Range("J16").Value = 11.59774691042095
$ws.Range("B17").Value = 21.93671533900496
$ws.Range("C17").Value = 14.26799382898495
$ws.Range("D17").Value = 6.123312757690691
$ws.Range("F17").Value = 58.93568026294059
$ws.Range("G17").Value = 3.701687248022912
$ws.Range("J17").Value = 11.58761245149184
$ws.Range("B18").Value = 21.83363226313837
$ws.Range("C18").Value = 14.17846545239699
$ws.Range("D18").Value = 6.135707529459411
$ws.Range("F18").Value = 58.67749614149249
$ws.Range("G18").Value = 3.702914436406917
$ws.Range("J18").Value = 11.58204592250648
$ws.Range("B19").Value = 21.79877099234557
$ws.Range("C19").Value = 14.14813677045615
$ws.Range("D19").Value = 6.139933342663905
$ws.Range("F19").Value = 58.59003252164158
$ws.Range("G19").Value = 3.703332543388482
$ws.Range("J19").Value = 11.58020630805476
$ws.Range("B20").Value = 21.95581211778396
$ws.Range("C20").Value = 14.2845550170636
$ws.Range("D20").Value = 6.12103260667541
$ws.Range("F20").Value = 58.98344036452578
$ws.Range("G20").Value = 3.701461357584603
$ws.Range("J20").Value = 11.58866411579395
$ws.Range("B21").Value = 22.48569128826786
$ws.Range("C21").Value = 14.7412349446236
$ws.Range("D21").Value = 6.059553209862503
$ws.Range("F21").Value = 60.30095765758227
$ws.Range("G21").Value = 3.695354586589036
$ws.Range("J21").Value = 11.62021579979582
$ws.Range("B22").Value = 22.83297643632115
$ws.Range("C22").Value = 15.03782866591265
$ws.Range("D22").Value = 6.02087752969646
$ws.Range("F22").Value = 61.15767400349593
$ws.Range("G22").Value = 3.691496722725625
$ws.Range("J22").Value = 11.64316703315904
$ws.Range("B23").Value = 22.64757790640358
$ws.Range("C23").Value = 14.87974139045595
$ws.Range("D23").Value = 6.041382050027968
$ws.Range("F23").Value = 60.70090055273016
$ws.Range("G23").Value = 3.69354360856121
$ws.Range("J23").Value = 11.63070474331752
$ws.Range("B24").Value = 21.94717793436622
$ws.Range("C24").Value = 14.2770681807172
$ws.Range("D24").Value = 6.122062917377407
$ws.Range("F24").Value = 58.96184934860407
$ws.Range("G24").Value = 3.7015634338469
$ws.Range("J24").Value = 11.58818784873925
$ws.Range("B25").Value = 21.20043134949723
$ws.Range("C25").Value = 13.62314695714889
$ws.Range("D25").Value = 6.21557884567738
$ws.Range("F25").Value = 57.07496484200021
$ws.Range("G25").Value = 3.710792824378815
$ws.Range("J25").Value = 11.55236411037291
